# Automatische test-sync: 2025-08-19 19:45:50
# Appends the new mail-log entry (row 9) to the "Logs" sheet, extends the
# conditional formatting ranges that cover the log table, and bumps the
# aggregate count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Logs sheet: add the new row --------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A9").Value = "CE-certificaten verzoek"
$logs.Range("B9").Value = "inkoop@testbedrijf123.nl"
$logs.Range("D9").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("F9").Value = "2025-08-19 19:45:42"
$logs.Range("G9").Value = "Nee"
$logs.Range("H9").Value = "Ja"
$logs.Range("I9").Value = "Nee"
$logs.Range("J9").Value = "Nee"

# --- 2. Extend the conditional-formatting ranges so row 9 is included ----
$ranges = @("D2:D8", "G2:G8", "H2:H8", "I2:I8", "J2:J8")
foreach ($oldRange in $ranges) {
    $newRange = $oldRange -replace "8$", "9"
    $fcs = $logs.Range($oldRange).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

# --- 3. Dashboard sheet: bump the aggregated count from 7 to 8 -----------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 8
